$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, pushing all existing data rows down by one.
$ws.Rows.Item(2).Insert()

# The inserted row inherits the (bold/border) formatting of the row above it;
# strip that so the new row starts out unstyled, like the other data rows.
$ws.Range("A2:E2").ClearFormats()

# Target data (row, date_serial, y0_year, y0_forecast, y1_year, y1_forecast)
$data = @(
  @(2, 39400, 2007, 1.75539628881467, 2008, 1.327368416067398),
  @(3, 39765, 2008, 2.213911448916162, 2009, 2.649257112350067),
  @(4, 40130, 2009, 2.533533936850563, 2010, 1.815660192323709),
  @(5, 40494, 2010, 2.088987486264915, 2011, 2.332261646026201),
  @(6, 40862, 2011, 1.212544822741002, 2012, 1.839804681163293),
  @(7, 41228, 2012, 1.196776590518644, 2013, 0.670590452940556),
  @(8, 41592, 2013, 0.4712609263772594, 2014, 0.8520644823059476),
  @(9, 41957, 2014, 0.8783377572271434, 2015, 1.474590898715178),
  @(10, 42321, 2015, 2.29066283401107, 2016, 2.597902967862775),
  @(11, 42689, 2016, 4.109890522944348, 2017, 3.628019428949036),
  @(12, 43053, 2017, 1.336316831462692, 2018, 1.626630409005325),
  @(13, 43418, 2018, 1.197912858979611, 2019, 1.216371234267344),
  @(14, 43783, 2019, 1.727537197898665, 2020, 2.164378481800822),
  @(15, 44159, 2020, 3.647228437274408, 2021, 3.845906281600109),
  @(16, 44525, 2021, 2.777797690741424, 2022, 1.875884305456199),
  @(17, 44890, 2022, 0.6994919452575576, 2023, -0.2388228654152447),
  @(18, 45254, 2023, -1.432689847121871, 2024, -0.7896638887521124),
  @(19, 45618, 2024, 2.033479419175133, 2025, 1.424898175306621)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
}

# Ensure the new row's date cell (A2) uses the same style as the other date cells
# (bold, centered, bordered, custom date format). The rest of row 2 (B2:E2) stays
# unstyled, matching the other plain data cells.
$ws.Range("A2").NumberFormat = $ws.Range("A3").NumberFormat
$ws.Range("A2").Font.Bold = $ws.Range("A3").Font.Bold
$ws.Range("A2").HorizontalAlignment = $ws.Range("A3").HorizontalAlignment
$ws.Range("A2").VerticalAlignment = $ws.Range("A3").VerticalAlignment
$ws.Range("A2").Borders.LineStyle = $ws.Range("A3").Borders.LineStyle
